$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 291
$ws1.Range("F3").Value = 1121
$ws1.Range("F4").Value = 2549
$ws1.Range("F5").Value = 219

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 291
$ws4.Range("F5").Value = 1121
$ws4.Range("F6").Value = 2549
$ws4.Range("F8").Value = 219
